$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.727.18'
$ws.Range('E2').Value = '  +0.77%  '
$ws.Range('D3').Value = '1.699.69'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.68'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4046'
$ws.Range('E8').Value = '  +1.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.511'
$ws.Range('E9').Value = '  -0.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.001'
$ws.Range('E10').Value = '  -0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.52'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08849'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.505'
$ws.Range('E13').Value = '  +4.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.64'
$ws.Range('E14').Value = '  +2.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.092'
$ws.Range('E15').Value = '  +7.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001322'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '1.696.66'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '99.31'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.79'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.090'
$ws.Range('E21').Value = '  +3.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.74'
$ws.Range('E23').Value = '  +5.22%  '
$ws.Range('D24').Value = '24.714.85'
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.148'
$ws.Range('E25').Value = '  +4.26%  '
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.69'
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.225'
$ws.Range('E28').Value = '  +23.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '164.28'
$ws.Range('E29').Value = '  +2.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '135.67'
$ws.Range('E30').Value = '  +1.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.151'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.054'
$ws.Range('E32').Value = '  +10.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09040'
$ws.Range('E33').Value = '  +6.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.072'
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02972'
$ws.Range('E35').Value = '  +8.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2775'
$ws.Range('E36').Value = '  +2.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.969'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('E38').Value = '  -1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.39'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09249'
$ws.Range('E40').Value = '  +2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.469'
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7758'
$ws.Range('E42').Value = '  +1.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '16.04'
$ws.Range('E43').Value = '  +5.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.7208'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.595'
$ws.Range('E45').Value = '  +3.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.213'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('B47').Value = 'Flow'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.355'
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('B48').Value = 'Frax'
$ws.Range('C48').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('E48').Value = '  +0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.78'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07987'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '89.66'
$ws.Range('E51').Value = '  +2.16%  '
